# Update 'Ofertas' worksheet: refresh iPhone 12 listings (new colors/prices/links),
# add 3 new listings (black iPhone 12, dupla-camera iPhone 12 re-link, RTX 3060 card re-link).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear all existing hyperlinks up front (cells will be re-populated + re-linked below).
$ws.Range("A1:C7").Hyperlinks.Delete()

$ws.Range("A2").Value = 'smartphone apple iphone 12, 64gb, branco, 5g, 6.1" super retina xdr oled, câmera dupla 12mp, selfie 12mp, ios 15'
$ws.Range("B2").Value = 3499
$ws.Range("C2").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABABGgJjZQ&sig=AOD64_34DxrFXzasoPQwqJBjphAkiWOBSg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8Iygw&adurl='
$ws.Hyperlinks.Add($ws.Range("C2"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABABGgJjZQ&sig=AOD64_34DxrFXzasoPQwqJBjphAkiWOBSg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8Iygw&adurl=')

$ws.Range("A3").Value = 'iphone 12 64gb roxo tela 6,1 4g câmera traseira 12mp vitrine'
$ws.Range("B3").Value = 3349
$ws.Range("C3").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAFGgJjZQ&sig=AOD64_06IUnEVYQQMWu4wJ0xAqMo7FPmJA&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8I0gw&adurl='
$ws.Hyperlinks.Add($ws.Range("C3"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAFGgJjZQ&sig=AOD64_06IUnEVYQQMWu4wJ0xAqMo7FPmJA&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8I0gw&adurl=')

$ws.Range("A4").Value = 'vitrine iphone 12 preto 64gb'
$ws.Range("B4").Value = 3199
$ws.Range("C4").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAVGgJjZQ&sig=AOD64_3yM2MV16Dy7D6qYOiDspoYT7HoHg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8I7gw&adurl='
$ws.Hyperlinks.Add($ws.Range("C4"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAVGgJjZQ&sig=AOD64_3yM2MV16Dy7D6qYOiDspoYT7HoHg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQww8I7gw&adurl=')

$ws.Range("A5").Value = 'iphone 12 64gb preto de vitrine tela 6,1 4g câmera traseira 12mp 12mp vitrine'
$ws.Range("B5").Value = 3349
$ws.Range("C5").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAnGgJjZQ&sig=AOD64_3va6pzH4P3F3x7UHlfSnxAZYVP0A&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I8hU&adurl='
$ws.Hyperlinks.Add($ws.Range("C5"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAnGgJjZQ&sig=AOD64_3va6pzH4P3F3x7UHlfSnxAZYVP0A&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I8hU&adurl=')

$ws.Range("A6").Value = 'iphone 12 64gb branco de vitrine tela 6,1&quot; 4g câmera traseira 12mp+12mp - vitrine'
$ws.Range("B6").Value = 3349
$ws.Range("C6").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABArGgJjZQ&sig=AOD64_13bJwCwSmGwViYbGTl_1F3GPe9EA&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I-RU&adurl='
$ws.Hyperlinks.Add($ws.Range("C6"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABArGgJjZQ&sig=AOD64_13bJwCwSmGwViYbGTl_1F3GPe9EA&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I-RU&adurl=')

$ws.Range("A7").Value = 'vitrine iphone 12 verde 64gb'
$ws.Range("B7").Value = 3199
$ws.Range("C7").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAtGgJjZQ&sig=AOD64_2-LB9KA2AzKKtTUoPfzPIrPYQsAg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I_BU&adurl='
$ws.Hyperlinks.Add($ws.Range("C7"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAtGgJjZQ&sig=AOD64_2-LB9KA2AzKKtTUoPfzPIrPYQsAg&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I_BU&adurl=')

$ws.Range("A8").Value = 'celular apple iphone 12 black 64gb vitrine/seminovo com carrregador e cabo'
$ws.Range("B8").Value = 3379
$ws.Range("C8").Value = 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAvGgJjZQ&sig=AOD64_3omzBLCWaDHRB9e48BnQ4AHFjw5w&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I_xU&adurl='
$ws.Hyperlinks.Add($ws.Range("C8"), 'https://www.google.com.br/aclk?sa=l&ai=DChcSEwj9isjRy8GAAxVcQkgAHfvrDVUYABAvGgJjZQ&sig=AOD64_3omzBLCWaDHRB9e48BnQ4AHFjw5w&ctype=5&q=&ved=0ahUKEwjQj8PRy8GAAxUMp5UCHVsiAtYQ9A4I_xU&adurl=')

$ws.Range("A9").Value = 'smartphone apple iphone 12 64gb câmera dupla'
$ws.Range("B9").Value = 3023
$ws.Range("C9").Value = 'https://www.buscape.com.br/celular/smartphone-apple-iphone-12-64gb-ios?_lc=88&searchterm=iphone%2012%2064%20gb'
$ws.Hyperlinks.Add($ws.Range("C9"), 'https://www.buscape.com.br/celular/smartphone-apple-iphone-12-64gb-ios?_lc=88&searchterm=iphone%2012%2064%20gb')

$ws.Range("A10").Value = 'placa de video nvidia geforce rtx 3060 ti 8 gb gddr6 192 bits asus dual-rtx3060ti-o8g-v2'
$ws.Range("B10").Value = 4108.27
$ws.Range("C10").Value = 'https://www.buscape.com.br/placa-de-video/placa-de-video-nvidia-geforce-rtx-3060-ti-8-gb-gddr6-192-bits-asus-dual-rtx3060ti-o8g-v2?_lc=88&searchterm=rtx%203060'
$ws.Hyperlinks.Add($ws.Range("C10"), 'https://www.buscape.com.br/placa-de-video/placa-de-video-nvidia-geforce-rtx-3060-ti-8-gb-gddr6-192-bits-asus-dual-rtx3060ti-o8g-v2?_lc=88&searchterm=rtx%203060')
